$wb = $excel.ActiveWorkbook

# --- Add the new, empty 'CaseDetailStat' sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsStat = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsStat.Name = "CaseDetailStat"

# --- Add the new 'CaseDetailStat_Message' sheet right after 'CaseDetailStat' ---
$wsMsg = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsStat)
$wsMsg.Name = "CaseDetailStat_Message"

# --- Populate column A of CaseDetailStat_Message with the connection/query log ---
# (three repeated Neo4j_URL/User_name/PWD/Cypher/Output blocks; the third block
#  reports a validation error with an empty Cypher query)
$wsMsg.Range("A1").Value = 'Neo4j_URL:'
$wsMsg.Range("A2").Value = 'bolt://ncidb-q325-c.nci.nih.gov:7687'
$wsMsg.Range("A3").Value = 'User_name:'
$wsMsg.Range("A4").Value = 'neo4j'
$wsMsg.Range("A5").Value = 'PWD:'
$wsMsg.Range("A6").Value = 'icdcDBneo4j0'
$wsMsg.Range("A7").Value = 'Cypher:'
$wsMsg.Range("A8").Value = 'MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN [''Adenocarcinoma of the cervix''] RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(t.clinical_trial_designation ,'''')as `Trial Code` , coalesce(a.arm_id,'''') As `Arm` , coalesce(a.arm_drug,'''') As `Arm Treatment` , coalesce(c.disease,'''') As Diagnosis , coalesce(c.gender,'''') As Gender , coalesce(c.race,'''') As Race , coalesce(c.ethnicity,'''') As Ethnicity'
$wsMsg.Range("A9").Value = 'Output:'
$wsMsg.Range("A10").Value = 'C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC01_Trials_Filter_Diagnosis-AdenoCervix_Neo4jData.xlsx'
$wsMsg.Range("A11").Value = 'Neo4j_URL:'
$wsMsg.Range("A12").Value = 'bolt://ncidb-q325-c.nci.nih.gov:7687'
$wsMsg.Range("A13").Value = 'User_name:'
$wsMsg.Range("A14").Value = 'neo4j'
$wsMsg.Range("A15").Value = 'PWD:'
$wsMsg.Range("A16").Value = 'icdcDBneo4j0'
$wsMsg.Range("A17").Value = 'Cypher:'
$wsMsg.Range("A18").Value = 'MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN [''Adenocarcinoma of the cervix''] OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial'
$wsMsg.Range("A19").Value = 'Output:'
$wsMsg.Range("A20").Value = 'C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC01_Trials_Filter_Diagnosis-AdenoCervix_Neo4jData.xlsx'
$wsMsg.Range("A21").Value = 'Cypher query should not be an empty string'
$wsMsg.Range("A22").Value = 'Neo4j_URL:'
$wsMsg.Range("A23").Value = 'bolt://ncidb-q325-c.nci.nih.gov:7687'
$wsMsg.Range("A24").Value = 'User_name:'
$wsMsg.Range("A25").Value = 'neo4j'
$wsMsg.Range("A26").Value = 'PWD:'
$wsMsg.Range("A27").Value = 'icdcDBneo4j0'
$wsMsg.Range("A28").Value = 'Cypher:'
$wsMsg.Range("A29").Formula = "=" + """" + """"
$wsMsg.Range("A30").Value = 'Output:'
$wsMsg.Range("A31").Value = 'C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC01_Trials_Filter_Diagnosis-AdenoCervix_Neo4jData.xlsx'
